# Refactor the Kyorugi/Poomsae entry-type tab: merge the 공인품새/자유품새
# "종목" distinction into the 참가부 labels, add 일반부 age groups, and append
# the new 개인전/복식전/단체전(자유품새) block (rows 20-31).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New 종목 (A) / 참가부 (B) pairs for rows 2-31, in order.
$rows = @(
    @("개인전", "초등부"),
    @("개인전", "중등부"),
    @("개인전", "고등부"),
    @("개인전", "대학부"),
    @("개인전", "일반부"),
    @("복식전", "초등부"),
    @("복식전", "중등부"),
    @("복식전", "고등부"),
    @("복식전", "대학부"),
    @("복식전", "일반부"),
    @("단체전", "초등부"),
    @("단체전", "중등부"),
    @("단체전", "고등부"),
    @("단체전", "대학부"),
    @("단체전", "일반부"),
    @("개인전(자유품새)", "초등부"),
    @("개인전(자유품새)", "중등부"),
    @("개인전(자유품새)", "고등부"),
    @("개인전(자유품새)", "대학부"),
    @("개인전(자유품새)", "일반부"),
    @("복식전(자유품새)", "초등부"),
    @("복식전(자유품새)", "중등부"),
    @("복식전(자유품새)", "고등부"),
    @("복식전(자유품새)", "대학부"),
    @("복식전(자유품새)", "일반부"),
    @("단체전(자유품새)", "초등부"),
    @("단체전(자유품새)", "중등부"),
    @("단체전(자유품새)", "고등부"),
    @("단체전(자유품새)", "대학부"),
    @("단체전(자유품새)", "일반부")
)

$r = 2
foreach ($pair in $rows) {
    $ws.Range("A$r").Value = $pair[0]
    $ws.Range("B$r").Value = $pair[1]
    # The old "세부부별" (C) values only made sense for the previous
    # 공인품새/자유품새 split; the merged layout no longer fills it in.
    $ws.Range("C$r").Value = ""
    $r = $r + 1
}
